$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    "Languages And Count",
    "Croatian(1)",
    "Azerbaijani(477)",
    "Dutch(2,126)",
    "Polish(2,133)",
    "Portuguese (Brazil)(2,389)",
    "Chinese (China)(2,347)",
    "Hindi(2,218)",
    "Swedish(2,127)",
    "Vietnamese(505)",
    "Hebrew(6)",
    "Kazakh(1,948)",
    "Ukrainian(2,156)",
    "Portuguese (Portugal)(1)",
    "French(2,402)",
    "Italian(2,160)",
    "Spanish(2,669)",
    "Chinese (Traditional)(45)",
    "Korean(2,235)",
    "English(5,859)",
    "Pushto(259)",
    "Catalan(1)",
    "Arabic(2,403)",
    "Uzbek(567)",
    "Russian(2,259)",
    "Turkish(2,177)",
    "Indonesian(2,313)",
    "Greek(2,127)",
    "Thai(2,140)",
    "German(2,331)",
    "Hungarian(1,372)",
    "Japanese(2,314)"
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 1).Value = $values[$i]
}

$ws.Range("C5").Select()
